{"js": "// Replace the Orion constellation viewing-window sentence in every\n// occurrence throughout the document body:\n//   \"V roku S\u00fahvezdie Orion 2022: ...\" ->\n//   \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Orion: ...\"\nconst oldText =\n  \"V roku S\u00fahvezdie Orion 2022: 16. \u2013 25. janu\u00e1ra, 14. \u2013 23. febru\u00e1ra, 14. \u2013 24. marca\";\nconst newText =\n  \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Orion: 16. \u2013 25. janu\u00e1ra, 14. \u2013 23. febru\u00e1ra, 14. \u2013 24. marca\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the Orion constellation viewing-window sentence everywhere it\n# appears in the document body:\n#   \"V roku S\u00fahvezdie Orion 2022: ...\" ->\n#   \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Orion: ...\"\n$d = $word.ActiveDocument\n\n$oldText = \"V roku S\u00fahvezdie Orion 2022: 16. \u2013 25. janu\u00e1ra, 14. \u2013 23. febru\u00e1ra, 14. \u2013 24. marca\"\n$newText = \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Orion: 16. \u2013 25. janu\u00e1ra, 14. \u2013 23. febru\u00e1ra, 14. \u2013 24. marca\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n\n$find.Execute(\n    [ref]$oldText,          # FindText\n    [ref]$false,            # MatchCase\n    [ref]$false,            # MatchWholeWord\n    [ref]$false,            # MatchWildcards\n    [ref]$false,            # MatchSoundsLike\n    [ref]$false,            # MatchAllWordForms\n    [ref]$true,             # Forward\n    [ref]1,                 # Wrap (wdFindContinue)\n    [ref]$false,            # Format\n    [ref]$newText,          # ReplaceWith\n    [ref]2                  # Replace (wdReplaceAll)\n) | Out-Null\n"}
